# Refresh cached Universalis market-price snapshots (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -- columns H:N)
# for the Leve rows whose prices moved since the last scheduled pull.
# One block per (sheet, Leve row); values taken from the refreshed feed.
$wb = $excel.ActiveWorkbook

### Sheet: ALC ###
$ws = $wb.Worksheets.Item("ALC")
# row 12
$ws.Range("H12").Value = 128.57143
$ws.Range("I12").Value = 120
$ws.Range("K12").Value = 120
$ws.Range("M12").Value = 50

$ws = $wb.Worksheets.Item("ALC")
# row 98
$ws.Range("H98").Value = 2993.111
$ws.Range("I98").Value = 867.75
$ws.Range("J98").Value = 6084.5454
$ws.Range("K98").Value = 867.75
$ws.Range("L98").Value = 6084.5454
$ws.Range("M98").Value = 630.25
$ws.Range("N98").Value = -9080.545399999999

$ws = $wb.Worksheets.Item("ALC")
# row 100
$ws.Range("H100").Value = 3450.25
$ws.Range("I100").Value = 3125.5
$ws.Range("J100").Value = 3775
$ws.Range("K100").Value = 3125.5
$ws.Range("L100").Value = 3775
$ws.Range("M100").Value = -2584.5
$ws.Range("N100").Value = -4857

$ws = $wb.Worksheets.Item("ALC")
# row 113
$ws.Range("H113").Value = 1893.6428
$ws.Range("I113").Value = 1744.375
$ws.Range("J113").Value = 2092.6667
$ws.Range("K113").Value = 1744.375
$ws.Range("L113").Value = 2092.6667
$ws.Range("M113").Value = 1509.625
$ws.Range("N113").Value = -8600.6667

$ws = $wb.Worksheets.Item("ALC")
# row 116
$ws.Range("H116").Value = 2332.9167
$ws.Range("I116").Value = 1828
$ws.Range("J116").Value = 2465.7896
$ws.Range("K116").Value = 1828
$ws.Range("L116").Value = 2465.7896
$ws.Range("M116").Value = 1614
$ws.Range("N116").Value = -9349.7896

$ws = $wb.Worksheets.Item("ALC")
# row 122
$ws.Range("H122").Value = 2993.111
$ws.Range("I122").Value = 867.75
$ws.Range("J122").Value = 6084.5454
$ws.Range("K122").Value = 2603.25
$ws.Range("L122").Value = 18253.6362
$ws.Range("M122").Value = -153.25
$ws.Range("N122").Value = -23153.6362

$ws = $wb.Worksheets.Item("ALC")
# row 132
$ws.Range("H132").Value = 4817512
$ws.Range("I132").Value = 7826526
$ws.Range("K132").Value = 23479578
$ws.Range("M132").Value = -23477048


### Sheet: ARM ###
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 994.73914
$ws.Range("I2").Value = 1006
$ws.Range("K2").Value = 1006
$ws.Range("M2").Value = -893

$ws = $wb.Worksheets.Item("ARM")
# row 97
$ws.Range("H97").Value = 1196.44
$ws.Range("I97").Value = 603.8095
$ws.Range("J97").Value = 4307.75
$ws.Range("K97").Value = 603.8095
$ws.Range("L97").Value = 4307.75
$ws.Range("M97").Value = -107.8095
$ws.Range("N97").Value = -5299.75

$ws = $wb.Worksheets.Item("ARM")
# row 116
$ws.Range("H116").Value = 994.73914
$ws.Range("I116").Value = 1006
$ws.Range("K116").Value = 1006
$ws.Range("M116").Value = 1288


### Sheet: BSM ###
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 994.73914
$ws.Range("I3").Value = 1006
$ws.Range("K3").Value = 1006
$ws.Range("M3").Value = -892

$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 1793.75
$ws.Range("J94").Value = 1957.1428
$ws.Range("L94").Value = 1957.1428
$ws.Range("N94").Value = -2859.1428

$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 1800
$ws.Range("I99").Value = 1968.75
$ws.Range("J99").Value = 1607.1428
$ws.Range("K99").Value = 1968.75
$ws.Range("L99").Value = 1607.1428
$ws.Range("M99").Value = -470.75
$ws.Range("N99").Value = -4603.1428


### Sheet: CRP ###
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 4463.8
$ws.Range("I31").Value = 5574.048
$ws.Range("J31").Value = 2798.4285
$ws.Range("K31").Value = 5574.048
$ws.Range("L31").Value = 2798.4285
$ws.Range("M31").Value = -5279.048
$ws.Range("N31").Value = -3388.4285

$ws = $wb.Worksheets.Item("CRP")
# row 34
$ws.Range("H34").Value = 4463.8
$ws.Range("I34").Value = 5574.048
$ws.Range("J34").Value = 2798.4285
$ws.Range("K34").Value = 5574.048
$ws.Range("L34").Value = 2798.4285
$ws.Range("M34").Value = -5372.048
$ws.Range("N34").Value = -3202.4285


### Sheet: CUL ###
$ws = $wb.Worksheets.Item("CUL")
# row 109
$ws.Range("H109").Value = 5201.9414
$ws.Range("I109").Value = 1287
$ws.Range("J109").Value = 6833.1665
$ws.Range("K109").Value = 3861
$ws.Range("L109").Value = 20499.4995
$ws.Range("M109").Value = -2821
$ws.Range("N109").Value = -22579.4995

$ws = $wb.Worksheets.Item("CUL")
# row 119
$ws.Range("H119").Value = 1414.1428
$ws.Range("I119").Value = 437.8
$ws.Range("K119").Value = 1313.4
$ws.Range("M119").Value = 3524.6


### Sheet: GSM ###
$ws = $wb.Worksheets.Item("GSM")
# row 113
$ws.Range("H113").Value = 2422.5833
$ws.Range("I113").Value = 2431
$ws.Range("J113").Value = 2330
$ws.Range("K113").Value = 2431
$ws.Range("L113").Value = 2330
$ws.Range("M113").Value = -261
$ws.Range("N113").Value = -6670


### Sheet: LTW ###
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 6166.6665
$ws.Range("I40").Value = 6166.6665
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6166.6665
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6030.6665
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# row 68
$ws.Range("H68").Value = 2964.5454
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2964.5454
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2964.5454
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4462.5454

$ws = $wb.Worksheets.Item("LTW")
# row 71
$ws.Range("H71").Value = 2964.5454
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2964.5454
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14822.727
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -22310.727

$ws = $wb.Worksheets.Item("LTW")
# row 136
$ws.Range("H136").Value = 3031.9778
$ws.Range("I136").Value = 2511.95
$ws.Range("J136").Value = 3448
$ws.Range("K136").Value = 7535.849999999999
$ws.Range("L136").Value = 10344
$ws.Range("M136").Value = -4985.849999999999
$ws.Range("N136").Value = -15444


### Sheet: WVR ###
$ws = $wb.Worksheets.Item("WVR")
# row 58
$ws.Range("H58").Value = 54347.668
$ws.Range("I58").Value = 10542.5
$ws.Range("J58").Value = 76250.25
$ws.Range("K58").Value = 10542.5
$ws.Range("L58").Value = 76250.25
$ws.Range("M58").Value = -10234.5
$ws.Range("N58").Value = -76866.25

$ws = $wb.Worksheets.Item("WVR")
# row 100
$ws.Range("H100").Value = 2082.3333
$ws.Range("I100").Value = 2428.9
$ws.Range("J100").Value = 349.5
$ws.Range("K100").Value = 4857.8
$ws.Range("L100").Value = 699
$ws.Range("M100").Value = -4316.8
$ws.Range("N100").Value = -1781

$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 588.4
$ws.Range("I113").Value = 497.5
$ws.Range("J113").Value = 724.75
$ws.Range("K113").Value = 1492.5
$ws.Range("L113").Value = 2174.25
$ws.Range("M113").Value = 677.5
$ws.Range("N113").Value = -6514.25

$ws = $wb.Worksheets.Item("WVR")
# row 136
$ws.Range("H136").Value = 5287.5
$ws.Range("I136").Value = 6157.0454
$ws.Range("J136").Value = 2099.1667
$ws.Range("K136").Value = 18471.1362
$ws.Range("L136").Value = 6297.500100000001
$ws.Range("M136").Value = -15921.1362
$ws.Range("N136").Value = -11397.5001
